$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Hydrogen ---
# D3: 161.8747113913984 -> blank (kept as an empty text cell, matching the
# sheet's existing convention of typed-but-empty cells). Using a
# quote-prefixed empty value yields an empty-string cell instead of fully
# clearing it, then the style is reset so no stray quote-prefix formatting
# is left behind.
$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = $ws.Range("C3").Style

# --- Row 4: Methanol ---
# C4: 16758.2131149072 -> 0
$ws.Range("C4").Value = 0

# --- Row 5: Ammonia ---
# C5: 68405.51980732256 -> 1922.932062252702
$ws.Range("C5").Value = 1922.932062252702

# --- Row 7: renamed Other -> Biogas, with a corrected value ---
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 273.8244887814103

# --- New row 8: Other (re-added below Biogas) ---
$ws.Range("A8").Value = "Other"
# Copy the header-cell formatting from A7 so A8 reuses the same bold/
# bordered/centered style instead of being left with default formatting.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)   # xlPasteFormats

# B8/C8 stay blank like the rest of the sheet's unused intersections, but
# keep the sheet's convention of an explicit empty-text cell rather than a
# fully absent one.
$ws.Range("B8").Value = "'"
$ws.Range("C8").Value = "'"
$ws.Range("B8").Style = $ws.Range("B7").Style
$ws.Range("C8").Style = $ws.Range("C7").Style

$ws.Range("D8").Value = 738.7221540173808
